$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue $ws "D2" "248.72"
Set-TextValue $ws "G2" "14"
Set-TextValue $ws "D3" "22.72"
Set-TextValue $ws "G3" "14"
Set-TextValue $ws "D4" "5.267"
Set-TextValue $ws "G4" "14"
Set-TextValue $ws "G5" "14"
Set-TextValue $ws "G6" "14"
Set-TextValue $ws "D7" "6.333"
Set-TextValue $ws "G7" "14"
Set-TextValue $ws "D8" "0.8053"
Set-TextValue $ws "G8" "14"
Set-TextValue $ws "D9" "0.9031"
Set-TextValue $ws "G9" "14"
Set-TextValue $ws "D10" "0.1411"
Set-TextValue $ws "G10" "14"
Set-TextValue $ws "D11" "0.07440"
Set-TextValue $ws "G11" "14"
Set-TextValue $ws "D12" "0.03094"
Set-TextValue $ws "G12" "14"
Set-TextValue $ws "D13" "0.03003"
Set-TextValue $ws "G13" "14"
Set-TextValue $ws "D14" "0.09380"
Set-TextValue $ws "G14" "14"
Set-TextValue $ws "D15" "3.861"
Set-TextValue $ws "G15" "14"
Set-TextValue $ws "D16" "0.001596"
Set-TextValue $ws "G16" "14"
Set-TextValue $ws "D17" "0.04771"
Set-TextValue $ws "G17" "14"
Set-TextValue $ws "G18" "14"
Set-TextValue $ws "D19" "0.0005800"
Set-TextValue $ws "G19" "14"
Set-TextValue $ws "D20" "0.006446"
Set-TextValue $ws "G20" "14"
Set-TextValue $ws "D21" "0.004991"
Set-TextValue $ws "G21" "14"
Set-TextValue $ws "D22" "0.0009992"
Set-TextValue $ws "G22" "14"
Set-TextValue $ws "G23" "14"
Set-TextValue $ws "D24" "3.695"
Set-TextValue $ws "G24" "14"
Set-TextValue $ws "D25" "2.201"
Set-TextValue $ws "G25" "14"
Set-TextValue $ws "D26" "0.3257"
Set-TextValue $ws "G26" "14"
Set-TextValue $ws "D27" "0.1292"
Set-TextValue $ws "G27" "14"
Set-TextValue $ws "G28" "14"
Set-TextValue $ws "G29" "14"
Set-TextValue $ws "G30" "14"
Set-TextValue $ws "G31" "14"
Set-TextValue $ws "G32" "14"
Set-TextValue $ws "G33" "14"
Set-TextValue $ws "G34" "14"
Set-TextValue $ws "G35" "14"
Set-TextValue $ws "G36" "14"
Set-TextValue $ws "G37" "14"
Set-TextValue $ws "G38" "14"
Set-TextValue $ws "G39" "14"
Set-TextValue $ws "D40" "0.03965"
Set-TextValue $ws "G40" "14"
Set-TextValue $ws "D41" "0.003041"
Set-TextValue $ws "E41" "40KickTokenKICKWorstin24h"
Set-TextValue $ws "G41" "14"
Set-TextValue $ws "G42" "14"
Set-TextValue $ws "D43" "0.002731"
Set-TextValue $ws "G43" "14"
Set-TextValue $ws "D44" "0.007709"
Set-TextValue $ws "G44" "14"
Set-TextValue $ws "D45" "0.00005583"
Set-TextValue $ws "G45" "14"
Set-TextValue $ws "G46" "14"
Set-TextValue $ws "D47" "0.4989"
Set-TextValue $ws "E47" "46CoinbaseStockTokenCOIN"
Set-TextValue $ws "G47" "14"
Set-TextValue $ws "D48" "0.2018"
Set-TextValue $ws "G48" "14"
Set-TextValue $ws "D49" "0.00002100"
Set-TextValue $ws "G49" "14"
Set-TextValue $ws "D50" "0.01010"
Set-TextValue $ws "G50" "14"
Set-TextValue $ws "G51" "14"
